$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.697917
$ws.Range("H2").Value = 2.093751
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8229576666666668
$ws.Range("N2").Value = 2.468873
$ws.Range("O2").Value = 0.03362764644735265
$ws.Range("P2").Value = 0.03362764644735265
$ws.Range("Q2").Value = 0.5743561458470001
$ws.Range("R2").Value = 5.169205312623001
$ws.Range("S2").Value = 0.03362764644735265
$ws.Range("T2").Value = 0.03362764644735265

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.697917
$ws.Range("H3").Value = 2.093751
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 23.10177466666667
$ws.Range("N3").Value = 69.305324
$ws.Range("O3").Value = 0.9439833204831614
$ws.Range("P3").Value = 0.9439833204831615
$ws.Range("Q3").Value = 16.123121270036
$ws.Range("R3").Value = 145.108091430324
$ws.Range("S3").Value = 0.9439833204831614
$ws.Range("T3").Value = 0.9439833204831615

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.697917
$ws.Range("H4").Value = 2.093751
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.5479189999999999
$ws.Range("N4").Value = 1.643757
$ws.Range("O4").Value = 0.02238903306948597
$ws.Range("P4").Value = 0.02238903306948598
$ws.Range("Q4").Value = 0.382401984723
$ws.Range("R4").Value = 3.441617862507
$ws.Range("S4").Value = 0.02238903306948597
$ws.Range("T4").Value = 0.02238903306948598
